$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each entry: row number, new "Periodo Mora" text (col E), new "Valor Mora" amount (col F)
# Column G ("Salario Basico") is updated uniformly to 781242 for all data rows.
$data = @(
    @(16, "1604", 6900),
    @(17, "1605", 6900),
    @(18, "1606", 6900),
    @(19, "1607", 6894),
    @(20, "1608", 6894),
    @(21, "1609", 6894),
    @(22, "1610", 6894),
    @(23, "1611", 6894),
    @(24, "1612", 6894),
    @(25, "1701", 13789),
    @(26, "1702", 13789),
    @(27, "1703", 13789),
    @(28, "1706", 14754),
    @(29, "1707", 14754),
    @(30, "1708", 14754),
    @(31, "1709", 14754),
    @(32, "1710", 14754),
    @(33, "1711", 14754),
    @(34, "1712", 14754),
    @(35, "1801", 22132),
    @(36, "1802", 22132),
    @(37, "1803", 22132),
    @(38, "1804", 22132),
    @(39, "1805", 22132),
    @(40, "1806", 22132),
    @(41, "1807", 22132),
    @(42, "1808", 22132),
    @(43, "1809", 23437),
    @(44, "1810", 23437),
    @(45, "1811", 23437),
    @(46, "1812", 23437),
    @(47, "1901", 31249),
    @(48, "1902", 31249),
    @(49, "1903", 31249),
    @(50, "1904", 31249),
    @(51, "1905", 31249),
    @(52, "1906", 31249),
    @(53, "1907", 31249),
    @(54, "1908", 31249),
    @(55, "1909", 31249),
    @(56, "1910", 31249),
    @(57, "1911", 31249),
    @(58, "1912", 31249),
    @(59, "2001", 31249),
    @(60, "2002", 31249),
    @(61, "2003", 31249)
)

foreach ($item in $data) {
    $row = $item[0]
    $periodo = $item[1]
    $valorMora = $item[2]

    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $valorMora
    $ws.Cells.Item($row, 7).Value = 781242
}
